$wb = $excel.ActiveWorkbook

$newNames = @(
    "summ41866441",
    "summ41966510",
    "summ42065146",
    "summ42165055",
    "summ42268301",
    "summ42368319",
    "summ42468129",
    "summ42669931",
    "summ42775797",
    "summ42875367",
    "summ42972025",
    "summ43069895",
    "summ43167770",
    "summ43268399",
    "summ43368289",
    "summ43468819",
    "summ43569714",
    "summ43668578",
    "summ43768770",
    "summ43875747",
    "summ43980833",
    "summ44082160",
    "summ44184440",
    "summ44293055",
    "summ44415828",
    "summ44538045",
    "summ44668703",
    "summ44802295",
    "summ44941118",
    "summ45083695",
    "summ45216239",
    "summ45349576",
    "summ45508893",
    "summ45664792",
    "summ45805362",
    "summ45942789",
    "summ46076815",
    "summ46203407",
    "summ46345005",
    "summ46505025",
    "summ46642384",
    "summ46783755",
    "summ46914276",
    "summ47055901",
    "summ47207239",
    "summ47347063",
    "summ47480729",
    "summ47614829",
    "summ47739365",
    "summ47858302"
)

for ($i = 1; $i -le $newNames.Count; $i++) {
    $wb.Worksheets.Item($i).Name = $newNames[$i - 1]
}
